$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as plain text so that
# numeric-looking values (e.g. "132.65", "1.00") are not coerced into
# real numbers (which would also lose fixed-point trailing zeros and
# introduce floating point artifacts).
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @{
    'D2' = '69.392.55'
    'E2' = '  +0.18%  '
    'D3' = '3.420.01'
    'E3' = '  +0.98%  '
    'E4' = '  +0.18%  '
    'D5' = '581.93'
    'E5' = '  -0.94%  '
    'D6' = '176.14'
    'E6' = '  -2.50%  '
    'E7' = '  +0.08%  '
    'D8' = '3.413.25'
    'E8' = '  +0.89%  '
    'E9' = '  -0.62%  '
    'D10' = '0.198'
    'E10' = '  +0.36%  '
    'E11' = '  -1.36%  '
    'D12' = '48.49'
    'E12' = '  -0.71%  '
    'E13' = '  -1.88%  '
    'D14' = '695.34'
    'E14' = '  +1.44%  '
    'D15' = '3.969.84'
    'E15' = '  +0.86%  '
    'D16' = '8.59'
    'E16' = '  -0.33%  '
    'D17' = '69.486.93'
    'E17' = '  +0.24%  '
    'D18' = '3.419.99'
    'E18' = '  +1.82%  '
    'E19' = '  +0.83%  '
    'D20' = '17.62'
    'E20' = '  -0.53%  '
    'E21' = '  -0.44%  '
    'E22' = '  -0.96%  '
    'D23' = '5.40'
    'E23' = '  -0.50%  '
    'D24' = '16.89'
    'E24' = '  -1.26%  '
    'D25' = '101.54'
    'E25' = '  -2.32%  '
    'E26' = '  -1.11%  '
    'E27' = '  -2.88%  '
    'D28' = '9.53'
    'E28' = '  -0.87%  '
    'D29' = '33.47'
    'E29' = '  -2.45%  '
    'D30' = '8.74'
    'E30' = '  +0.31%  '
    'D31' = '7.05'
    'E31' = '  +1.02%  '
    'D32' = '574.58'
    'E32' = '  +3.47%  '
    'E33' = '  -0.13%  '
    'D34' = '10.98'
    'E34' = '  -1.86%  '
    'D35' = '58.40'
    'E35' = '  +0.42%  '
    'E36' = '  -3.19%  '
    'E37' = '  +0.04%  '
    'D38' = '3.535.09'
    'E38' = '  -4.50%  '
    'E39' = '  -1.78%  '
    'D40' = '34.85'
    'E40' = '  -0.38%  '
    'D41' = '0.0₃0728'
    'E41' = '  +3.48%  '
    'E42' = '  +0.33%  '
    'D43' = '2.66'
    'E43' = '  -1.19%  '
    'B44' = 'ApeXProtocol'
    'C44' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D44' = '3.38'
    'E44' = '  +3.69%  '
    'B45' = 'TheGraph'
    'C45' = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    'D45' = '0.331'
    'E45' = '  -2.49%  '
    'B46' = 'VeChain'
    'C46' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D46' = '0.0416'
    'E46' = '  -1.71%  '
    'B47' = 'Mantle'
    'C47' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D47' = '1.44'
    'E47' = '  +4.22%  '
    'B48' = 'ThetaToken'
    'C48' = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
    'D48' = '2.63'
    'E48' = '  -0.56%  '
    'B49' = 'Stellar'
    'C49' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D49' = '0.128'
    'E49' = '  -1.29%  '
    'B50' = 'FirstDigitalUSD'
    'C50' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D50' = '1.00'
    'E50' = '  -0.14%  '
    'B51' = 'Monero'
    'C51' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D51' = '132.65'
    'E51' = '  +0.12%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

